$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row - copy formatting from H1 (bold, bordered, centered style)
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("I1").Value = "I0"

$ws.Range("H1").Copy()
$ws.Range("J1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("J1").Value = "IF"

# Data rows 2..38: I = 1 (constant), J = copy of H
for ($r = 2; $r -le 38; $r++) {
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $ws.Cells.Item($r, 8).Value()
}
